$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-number labels in column A (rows 2-5), keeping them as text
# so leading zeros are preserved.
$ws.Range("A2").Value = "000"
$ws.Range("A3").Value = "001"
$ws.Range("A4").Value = "002"
$ws.Range("A5").Value = "003"

# Move the active selection to A6
$ws.Range("A6").Select()
